$wb = $excel.ActiveWorkbook

# --- Sheet3 (Feuil3): insert a new column A (row numbers / letters) and
#     a new column B (short codes), shifting the old column A data to column C.
$ws3 = $wb.Worksheets.Item("Feuil3")

# Existing values in column A (before the edit) for rows 1-10, 12-14
$oldValues = @{
    1  = "Ligne"
    2  = "Modèle"
    3  = "Code pièce"
    4  = "Couleur"
    5  = "Taille"
    6  = "Quantité"
    7  = "Unité"
    8  = "Prix unitaire"
    9  = "Montant"
    10 = "Libellé"
    12 = "Total facture"
    13 = "Total montant"
    14 = "Référence facture"
}

# New column A values (numbers for rows 1-10, letters for rows 12-14)
$colA = @{
    1  = 1
    2  = 2
    3  = 3
    4  = 4
    5  = 5
    6  = 6
    7  = 7
    8  = 8
    9  = 9
    10 = 10
    12 = "A"
    13 = "B"
    14 = "C"
}

# New column B values (short codes)
$colB = @{
    1  = "LI"
    2  = "MO"
    3  = "PI"
    4  = "CO"
    5  = "TA"
    6  = "QT"
    7  = "UN"
    8  = "PU"
    9  = "MT"
    10 = "LI"
    12 = "TF"
    13 = "TM"
    14 = "RF"
}

foreach ($r in $colA.Keys) {
    $ws3.Range("C$r").Value = $oldValues[$r]
    $ws3.Range("B$r").Value = $colB[$r]
    $ws3.Range("A$r").Value = $colA[$r]
}

# Select A14 as before and mark this sheet as the active/selected tab.
$ws3.Range("A14").Select()
$ws3.Activate()

# --- Workbook-level defined name: tabCodes now points to column C (the
#     textual codes) instead of column A (which now holds row numbers).
$wb.Names.Item("tabCodes").RefersTo = "=Feuil3!`$C:`$C"

# --- Sheet2 (Feuil2): the shared formula in M20 no longer spans the whole
#     M20:M57 range - it is now split, with M20:M22 keeping si=12, and a new
#     shared formula (M24:M57) taking over the rest (M23 already had its own
#     standalone formula). Re-entering the formula in M20:M22 as an array of
#     individual formulas achieves the narrower shared-formula range, and
#     M24:M57 is re-entered so the existing shared group is rebuilt cleanly.
$ws2 = $wb.Worksheets.Item("Feuil2")

$ws2.Range("M20:M22").FormulaR1C1 = '=RIGHT("000" & TRIM(INDEX(INDIRECT(RC7),OFFSET(INDIRECT(RC5), 0, COLUMN()-1)+RC2-1)), 3)'
$ws2.Range("M24:M57").FormulaR1C1 = '=RIGHT("000" & TRIM(INDEX(INDIRECT(RC7),OFFSET(INDIRECT(RC5), 0, COLUMN()-1)+RC2-1)), 3)'

$wb.Save()
